$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020
